$wb = $excel.ActiveWorkbook

# "out" sheet gets a new row of data
$wsOut = $wb.Worksheets.Item("out")

$wsOut.Cells.Item(11, 1).Value = 20191010
$wsOut.Cells.Item(11, 2).Value = "ZL Asica"
$wsOut.Cells.Item(11, 3).Value = 80.2
$wsOut.Cells.Item(11, 4).Value = "服务器10月费用"

# Select the newly entered row and activate the "out" sheet/tab
$wsOut.Activate()
$wsOut.Range("D11").Select()
